$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new data row at row 11 (weekly update, pushes the
# existing Oct-04 "Primera" row and everything after it down by one) ---
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = "2021-10-04"
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100108
$ws.Range("H11").Value = "Tropicales y subtropicales"
$ws.Range("I11").Value = 100108003
$ws.Range("J11").Value = "Maracuyá"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 75000
$ws.Range("O11").Value = 75000
$ws.Range("P11").Value = 75000
$ws.Range("Q11").Value = "$/caja 18 kilos"
$ws.Range("R11").Value = "Perú"
$ws.Range("S11").Value = 4167
$ws.Range("T11").Value = 18

# --- Step 2: insert two new data rows at row 39 (new Nov-22 week, origin
# Región de Arica y Parinacota), pushing the remaining historical rows down ---
$ws.Range("A39:A40").EntireRow.Insert()

$ws.Range("A39").Value = 9
$ws.Range("B39").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C39").Value = "Metropolitana"
$ws.Range("D39").Value = "2021-11-22"
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100108
$ws.Range("H39").Value = "Tropicales y subtropicales"
$ws.Range("I39").Value = 100108003
$ws.Range("J39").Value = "Maracuyá"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 25
$ws.Range("N39").Value = 60000
$ws.Range("O39").Value = 60000
$ws.Range("P39").Value = 60000
$ws.Range("Q39").Value = "$/caja 18 kilos"
$ws.Range("R39").Value = "Región de Arica y Parinacota"
$ws.Range("S39").Value = 3333
$ws.Range("T39").Value = 18

$ws.Range("A40").Value = 9
$ws.Range("B40").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = "2021-11-22"
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100108
$ws.Range("H40").Value = "Tropicales y subtropicales"
$ws.Range("I40").Value = 100108003
$ws.Range("J40").Value = "Maracuyá"
$ws.Range("K40").Value = "Sin especificar"
$ws.Range("L40").Value = "Segunda"
$ws.Range("M40").Value = 20
$ws.Range("N40").Value = 56000
$ws.Range("O40").Value = 56000
$ws.Range("P40").Value = 56000
$ws.Range("Q40").Value = "$/caja 18 kilos"
$ws.Range("R40").Value = "Región de Arica y Parinacota"
$ws.Range("S40").Value = 3111
$ws.Range("T40").Value = 18
